$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "vijay" to "Sheet1"
$ws.Name = "Sheet1"

# Add new header "Login" in E1, reusing the same formatting as the other
# header cells (copy format from D1, then set the value/text).
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Login"

# E2 stays blank (empty string) for the existing row
$ws.Range("E2").Value = ""

# Add row 3, replicating row 2's match data with a login value in E3
$ws.Range("A3").Value = "22-03-2025"
$ws.Range("B3").Value = "Kolkata Knight Riders vs Royal Challengers Bengaluru"
$ws.Range("C3").Value = "Kolkata Knight Riders"
$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "vijay"
